$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "department" column (C) previously held the single value
# "FACULTY OF HOSPITALITY" for every promo row. Split it into the new
# per-row department labels used going forward.
$ws.Range("C2").Value = "Hospitality"
$ws.Range("C3").Value = "Packages"
$ws.Range("C4").Value = "Packages"
$ws.Range("C5").Value = "Packages"
